# Word COM-interop script applying the "Alt-Right" document edits described
# by the commit "added more to Altright".

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $found = $d.Content.Find.Execute(
        $find, $true, $false, $false, $false, $false,
        $true, 1, $false, $replace, 2)
    if (-not $found) {
        Write-Output ("WARNING: could not find -> " + $find)
    }
}

# --- Paragraph: "One of the more extremist groups... home on messageboards..." ---
Replace-Text "has its home on messageboards like" "has its home on message boards like"
Replace-Text "dangerous ideology because of it’s hateful rhetoric" "dangerous ideology because of its hateful rhetoric"

# --- Paragraph: "Many stories of Universities..." ---
Replace-Text "renaming buildings or citys taking down" "renaming buildings or cities taking down"
Replace-Text "intentionally misgendering a trans person" "intentionally mis-gendering a trans person"
Replace-Text "stories of violent muslim immigr" "stories of violent Muslim immigr"
Replace-Text "deserve an ethnostate because both Isreal and Japan" "deserve an ethno-state because both Israel and Japan"

# --- Paragraph: "As an onlooker to the birth and rise..." ---
Replace-Text "got out of hand very quicky." "got out of hand very quickly."
Replace-Text "much traction in it’s early years, until it took hold" "much traction in its early years until it took hold"
Replace-Text "such as online messageboards like 4chan and reddit" "such as online message boards like 4chan and reddit"
Replace-Text "evolution of the mindframe of the Alt-Right" "evolution of the mind frame of the Alt-Right"
Replace-Text "even the most henouis things" "even the most heinous things"
Replace-Text "comic by the name of  “Boy’s Club” by Matt Furie." "comic by the name of “Boy’s Club” by Matt Furie."

# --- Move the "_GoBack" bookmark from the end of the document to right after
#     "One o" in the opening sentence (mirrors Word re-marking the last edit
#     location after the author typed in that spot). ---
$findRange = $d.Content
$findRange.Find.Execute("One o")
$findRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $findRange)

# --- Append two new Heading1-styled paragraphs at the end of the document. ---
$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newPara1 = $d.Paragraphs.Last
$newRange1 = $newPara1.Range
$newRange1.Collapse(0)
$newRange1.InsertAfter("`tThe problem at the root of the Alt-Right movement is that it’s based upon the idea of “saving western (and by extension white) culture” and preserving it forever (a thousand year reich). This idea of preserving a culture to stay the same is a futile concept, especially in America. Youtuber Contrapoints said it best, “It’s like a kind of Egyptian mummification urge, the feeling that if we could just embalm our own culture and protect it from the whims of time, then it’ll never die. But just like mummifying dead bodies for the afterlife, this is a delusion” (Contrapoints). One of the main problems of preserving white western culture, is that as a culture, western culture is one of the most diverse, indulging in food from around the world, and most of the music being based off rock and blues, both being black inventions. Whiteness as a concept is flimsy at most, with in history people such as the Irish immigrants to the US would be considered white, and accepted today, just 150 years ago were discriminated against just as the white supremacists would discriminate against black people now. There is no doubt that the Alt-Right is rooted in white supremacism, especially when one looks at its supporters online, and their rhetoric. The ADL is also cited as saying “In fact, Alt Righters reject modern conservatism explicitly because they believe that mainstream conservatives are not advocating for the interests of white people as a group” (ADL).")

$newPara1 = $d.Paragraphs.Last
$tailRange = $newPara1.Range
$tailRange.Collapse(0)
$tailRange.InsertParagraphAfter()

$newPara2 = $d.Paragraphs.Last
$newRange2 = $newPara2.Range
$newRange2.Collapse(0)
$newRange2.InsertAfter("`t")

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
Write-Output "Done."
